# Generate Report for handoff
# Updates the status of the "9fc327a9-32f4-4496-87fb-22743577c874.md" file
# from "Handed back: in sync with en-US" to "Ready for handoff" on the
# Overview sheet as well as the per-locale sheets, and records the new
# handoff datetime on each locale sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to the 9fc327a9...md file
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: update status + latest handoff datetime for row 3
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-01-18 12:24:43"

# de-de sheet: update status + latest handoff datetime for row 3
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-01-18 12:24:52"
